$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.916.71'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.734.12'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '245.79'
$ws.Range('E5').Value = '  +3.04%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.5037'
$ws.Range('E7').Value = '  -2.53%  '
$ws.Range('D8').Value = '0.2734'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '0.06181'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').Value = '1.740.93'
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('D11').Value = '0.07245'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').Value = '0.6542'
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('D13').Value = '15.20'
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').Value = '4.784'
$ws.Range('E14').Value = '  +4.15%  '
$ws.Range('D15').Value = '77.15'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '0.9992'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '0.9988'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '25.932.65'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').Value = '11.93'
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').Value = '0.000006840'
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('D21').Value = '4.595'
$ws.Range('E21').Value = '  +7.59%  '
$ws.Range('D22').Value = '1.962.29'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = '8.816'
$ws.Range('E23').Value = '  +1.61%  '
$ws.Range('D24').Value = '5.495'
$ws.Range('E24').Value = '  +4.75%  '
$ws.Range('D25').Value = '134.04'
$ws.Range('E25').Value = '  -3.32%  '
$ws.Range('D26').Value = '15.27'
$ws.Range('E26').Value = '  +0.91%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '1.798'
$ws.Range('E27').Value = '  +1.96%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.424'
$ws.Range('E28').Value = '  -5.59%  '
$ws.Range('D29').Value = '105.73'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '3.988'
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('D31').Value = '0.08143'
$ws.Range('E31').Value = '  -1.84%  '
$ws.Range('D32').Value = '3.702'
$ws.Range('E32').Value = '  +1.53%  '
$ws.Range('D33').Value = '0.04750'
$ws.Range('E33').Value = '  +3.56%  '
$ws.Range('D34').Value = '2.653'
$ws.Range('D35').Value = '0.9998'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').Value = '0.6158'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').Value = '2.747'
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01613'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '0.8797'
$ws.Range('E39').Value = '  +18.96%  '
$ws.Range('D40').Value = '1.963'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('D41').Value = '0.9989'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').Value = '102.07'
$ws.Range('E42').Value = '  +4.29%  '
$ws.Range('D43').Value = '0.3925'
$ws.Range('E43').Value = '  +2.30%  '
$ws.Range('D44').Value = '5.025'
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').Value = '0.1182'
$ws.Range('E45').Value = '  +5.09%  '
$ws.Range('D46').Value = '6.388'
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').Value = '55.86'
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('D48').Value = '0.05284'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').Value = '30.86'
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('D50').Value = '0.3502'
$ws.Range('E50').Value = '  +2.82%  '
$ws.Range('D51').Value = '7.634'
$ws.Range('E51').Value = '  -0.02%  '
